# Add "フィブリノーゲン" (Fibrinogen) test row and a "■血液型" (Blood type)
# section (ABO blood type / Rh blood type, the latter with a +/- dropdown)
# to the lab-test template, inserted right after the "D-dimer" row (old
# row 24) and before the "■臨床化学検査" section (old row 25 -> new row 30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert five new blank rows at row 25; everything from the old row 25
# onward (the "■臨床化学検査" section and below) shifts down to row 30+.
$ws.Rows("25:29").Insert()

# --- Fix up formatting on the freshly inserted rows -------------------
# Row 25 pattern (A: label, B: blank value, C: unit) -> copy from row 3
# (Hb / g/dL), which already uses styles 6/7/8.
$ws.Range("A3:C3").Copy()
$ws.Range("A25:C25").PasteSpecial(-4122)

# Rows 26-28 pattern (A: label or blank, B: blank, C: blank) -> copy from
# row 31 ("■臨床化学検査" after the shift), styles 6/7/9.
$ws.Range("A31:C31").Copy()
$ws.Range("A26:C28").PasteSpecial(-4122)

# Row 29 pattern: A uses the "spacer" style (10) but still carries a
# label, B/C stay blank (7/9) -> copy A from row 18 (blank spacer) and
# B/C from row 31.
$ws.Range("A18").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("B31:C31").Copy()
$ws.Range("B29:C29").PasteSpecial(-4122)

# Row 30 pattern: full spacer row (A/B/C = 10/7/9) -> copy from row 18.
$ws.Range("A18:C18").Copy()
$ws.Range("A30:C30").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Restore the row heights lost by Insert().
$ws.Rows("25:30").RowHeight = 22

# --- Values -------------------------------------------------------------
$ws.Range("A25").Value = "フィブリノーゲン"
$ws.Range("C25").Value = "mg/dL"

$ws.Range("A27").Value = "■血液型"
$ws.Range("A28").Value = "ABO式血液型"
$ws.Range("A29").Value = "Rh式血液型"

# --- Data validation: +/- dropdown for Rh blood type value (B29) -------
$ws.Range("B29").Validation.Add(3, 1, 1, '"＋,−"')

Write-Output "done"
